$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9488581418991089
$ws.Range("B1").Value = 2.170935392379761
$ws.Range("C1").Value = 8.472025871276855
$ws.Range("D1").Value = 1.741669654846191
$ws.Range("E1").Value = 1.400412201881409
